# Hortaliza, Vega Modelo de Temuco - Cebollín
# Weekly data refresh: insert two new observation rows (465, 466) ahead of
# the existing history, pushing the previously-recorded rows down by two
# (old row 465 -> new row 467, ... old row 543 -> new row 545).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 465-466; Excel shifts everything from the old
# row 465 downward by two rows (465->467 ... 543->545), exactly matching
# the target layout.
$ws.Range("A465:A466").EntireRow.Insert()

# --- Populate new row 465 ---
$ws.Range("A465").Value2 = 10
$ws.Range("B465").Value2 = "Vega Modelo de Temuco"
$ws.Range("C465").Value2 = "La Araucanía"
$ws.Range("D465").Value2 = 44951
$ws.Range("E465").Value2 = 9
$ws.Range("F465").Value2 = 100112037
$ws.Range("G465").Value2 = "Cebollín"
$ws.Range("H465").Value2 = "Sin especificar"
$ws.Range("I465").Value2 = "Primera"
$ws.Range("J465").Value2 = 125
$ws.Range("K465").Value2 = 6000
$ws.Range("L465").Value2 = 6000
$ws.Range("M465").Value2 = 6000
$ws.Range("N465").Value2 = "`$/docena de paquetes"
$ws.Range("O465").Value2 = "Provincia de Cautín"
$ws.Range("P465").Value2 = 500
$ws.Range("Q465").Value2 = 12
$ws.Range("R465").Value2 = "Hortaliza"

# --- Populate new row 466 ---
$ws.Range("A466").Value2 = 10
$ws.Range("B466").Value2 = "Vega Modelo de Temuco"
$ws.Range("C466").Value2 = "La Araucanía"
$ws.Range("D466").Value2 = 44951
$ws.Range("E466").Value2 = 9
$ws.Range("F466").Value2 = 100112037
$ws.Range("G466").Value2 = "Cebollín"
$ws.Range("H466").Value2 = "Sin especificar"
$ws.Range("I466").Value2 = "Primera"
$ws.Range("J466").Value2 = 115
$ws.Range("K466").Value2 = 6000
$ws.Range("L466").Value2 = 6000
$ws.Range("M466").Value2 = 6000
$ws.Range("N466").Value2 = "`$/docena de paquetes"
$ws.Range("O466").Value2 = "Región de O'Higgins"
$ws.Range("P466").Value2 = 500
$ws.Range("Q466").Value2 = 12
$ws.Range("R466").Value2 = "Hortaliza"
